$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 681 (weekly update: new price observation),
# pushing the existing rows 681:782 down to 682:783.
$ws.Rows("681:681").Insert()

$ws.Range("A681").Value = 6
$ws.Range("B681").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C681").Value = 'Metropolitana'
$ws.Range("D681").Value = 45154
$ws.Range("E681").Value = 13
$ws.Range("F681").Value = 100112052
$ws.Range("G681").Value = 'Albahaca'
$ws.Range("H681").Value = 'Sin especificar'
$ws.Range("I681").Value = 'Primera'
$ws.Range("J681").Value = 38
$ws.Range("K681").Value = 4500
$ws.Range("L681").Value = 5000
$ws.Range("M681").Value = 4697
$ws.Range("N681").Value = '$/paquete'
$ws.Range("O681").Value = 'Región de Arica y Parinacota'
$ws.Range("P681").Value = 4697
$ws.Range("Q681").Value = 1
$ws.Range("R681").Value = 'Hortaliza'
